{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is: async (context) => { ... }\n\n// The document is a title paragraph (\"2025-04-15 Tuesday\") followed by a\n// single 20-row x 5-column table. Only 5 of the rows carry multiplication\n// problems (rows 0, 4, 9, 14, 19); the rest are spacer rows. We address\n// each cell by its (row, col) grid position so that the update is\n// unambiguous even though some new values duplicate other old values\n// elsewhere in the table (e.g. \"99\u00d789=\" both disappears and reappears).\nconst tableRowUpdates = [\n  { row: 0, values: [\"39\u00d713=\", \"63\u00d790=\", \"91\u00d725=\", \"89\u00d769=\", \"87\u00d769=\"] },\n  { row: 4, values: [\"72\u00d788=\", \"85\u00d770=\", \"98\u00d752=\", \"18\u00d711=\", \"61\u00d765=\"] },\n  { row: 9, values: [\"62\u00d759=\", \"52\u00d788=\", \"48\u00d791=\", \"48\u00d713=\", \"99\u00d789=\"] },\n  { row: 14, values: [\"90\u00d722=\", \"65\u00d729=\", \"51\u00d766=\", \"73\u00d774=\", \"52\u00d716=\"] },\n  { row: 19, values: [\"73\u00d733=\", \"90\u00d786=\", \"66\u00d764=\", \"25\u00d726=\", \"39\u00d731=\"] },\n];\n\nconst body = context.document.body;\n\n// Update the date line in the title paragraph.\nconst dateResults = body.search(\"2025-04-15 Tuesday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nfor (const range of dateResults.items) {\n  range.insertText(\"2025-04-22 Tuesday\", \"Replace\");\n}\nawait context.sync();\n\n// Update the multiplication table cells by explicit grid position.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const { row, values } of tableRowUpdates) {\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = values[col];\n  }\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n\n$d = $word.ActiveDocument\n\n# --- Update the date line in the title paragraph -------------------------\n$dateRange = $d.Content\n$found = $dateRange.Find.Execute(\n    \"2025-04-15 Tuesday\", $false, $false, $false, $false, $false,\n    $true, 1, $false, \"2025-04-22 Tuesday\", 2)\n\n# --- Update the multiplication table -------------------------------------\n# The document holds a single 20-row x 5-column table. Only rows 1, 5, 10,\n# 15 and 20 (1-based COM indices) carry multiplication problems; the rest\n# are blank spacer rows. Addressing each cell by its (row, col) position\n# keeps the update unambiguous even though some new values duplicate other\n# old values elsewhere in the table (e.g. \"99\u00d789=\" disappears from one\n# cell and reappears in another).\n\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @{\n    1  = @(\"39\u00d713=\", \"63\u00d790=\", \"91\u00d725=\", \"89\u00d769=\", \"87\u00d769=\")\n    5  = @(\"72\u00d788=\", \"85\u00d770=\", \"98\u00d752=\", \"18\u00d711=\", \"61\u00d765=\")\n    10 = @(\"62\u00d759=\", \"52\u00d788=\", \"48\u00d791=\", \"48\u00d713=\", \"99\u00d789=\")\n    15 = @(\"90\u00d722=\", \"65\u00d729=\", \"51\u00d766=\", \"73\u00d774=\", \"52\u00d716=\")\n    20 = @(\"73\u00d733=\", \"90\u00d786=\", \"66\u00d764=\", \"25\u00d726=\", \"39\u00d731=\")\n}\n\nforeach ($row in $rowUpdates.Keys) {\n    $values = $rowUpdates[$row]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n\n$d.Save()\n"}
